$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values so Excel keeps them as text
$textCells = @("D5", "D6", "D10", "D16", "D17", "D19", "D20", "D21", "D25", "D29", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D49", "D50")
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '67.648.94'
$ws.Range('E2').Value = '  -2.53%  '
$ws.Range('D3').Value = '2.426.15'
$ws.Range('E3').Value = '  -2.55%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '548.65'
$ws.Range('E5').Value = '  -3.00%  '
$ws.Range('D6').Value = '158.72'
$ws.Range('E6').Value = '  -3.04%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -3.13%  '
$ws.Range('D9').Value = '2.422.03'
$ws.Range('E9').Value = '  -2.65%  '
$ws.Range('D10').Value = '0.143'
$ws.Range('E10').Value = '  -9.02%  '
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('E12').Value = '  -6.12%  '
$ws.Range('E13').Value = '  -4.08%  '
$ws.Range('D14').Value = '2.869.88'
$ws.Range('E14').Value = '  -2.53%  '
$ws.Range('D15').Value = '67.503.66'
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('D16').Value = '0.0000163'
$ws.Range('E16').Value = '  -6.59%  '
$ws.Range('D17').Value = '22.88'
$ws.Range('E17').Value = '  -5.59%  '
$ws.Range('D18').Value = '2.421.28'
$ws.Range('E18').Value = '  -2.64%  '
$ws.Range('D19').Value = '10.64'
$ws.Range('E19').Value = '  -4.82%  '
$ws.Range('D20').Value = '336.19'
$ws.Range('E20').Value = '  -2.33%  '
$ws.Range('D21').Value = '6.90'
$ws.Range('E21').Value = '  -6.03%  '
$ws.Range('E22').Value = '  -3.89%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('E24').Value = '  -5.44%  '
$ws.Range('D25').Value = '65.88'
$ws.Range('E25').Value = '  -5.14%  '
$ws.Range('E26').Value = '  -7.22%  '
$ws.Range('D27').Value = '2.549.20'
$ws.Range('E27').Value = '  -2.56%  '
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('D29').Value = '7.93'
$ws.Range('E29').Value = '  -8.23%  '
$ws.Range('D30').Value = '0.0₃0801'
$ws.Range('E30').Value = '  -8.03%  '
$ws.Range('E31').Value = '  -9.14%  '
$ws.Range('D32').Value = '0.998'
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('D33').Value = '418.62'
$ws.Range('E33').Value = '  -5.17%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '1.11'
$ws.Range('E34').Value = '  -6.28%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').Value = '1.59'
$ws.Range('E35').Value = '  -6.20%  '
$ws.Range('D36').Value = '156.68'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('D37').Value = '18.97'
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('E39').Value = '  -5.08%  '
$ws.Range('D40').Value = '17.49'
$ws.Range('E40').Value = '  -3.45%  '
$ws.Range('D41').Value = '0.297'
$ws.Range('E41').Value = '  -5.39%  '
$ws.Range('D42').Value = '4.26'
$ws.Range('E42').Value = '  -6.96%  '
$ws.Range('D43').Value = '1.42'
$ws.Range('E43').Value = '  -10.27%  '
$ws.Range('D44').Value = '1.06'
$ws.Range('E44').Value = '  -0.58%  '
$ws.Range('D45').Value = '131.85'
$ws.Range('E46').Value = '  -7.83%  '
$ws.Range('E47').Value = '  -4.88%  '
$ws.Range('E48').Value = '  -2.79%  '
$ws.Range('D49').Value = '0.469'
$ws.Range('E49').Value = '  -8.24%  '
$ws.Range('D50').Value = '0.548'
$ws.Range('E50').Value = '  -3.88%  '
$ws.Range('E51').Value = '  -2.44%  '
